# Gain resistors updated to 1k, 100
#
# - R11 (previously its own "5k" group) is merged into the "1k" group
#   (designator list for row 7 gains "R11" in sorted position).
# - The "200" resistor group (row 10: R10, R19, R26) becomes a plain
#   numeric "100" value (left-aligned to match the surrounding text
#   cells visually).
# - The now-empty "5k" / R11 row (row 19) is removed entirely, shifting
#   all subsequent rows up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fold R11 into the 1k designator list (row 7, column B).
$ws.Range("B7").Value = "R1, R6, R7, R8, R9, R11, R12, R13, R18, R20, R21, R22, R23, R24, R27, R28, R29, R30, R32, R33, R36, R38, R40"

# 2) Change the "200" comment (row 10, column A) to the numeric value 100,
#    keeping it left-aligned like the rest of column A.
$ws.Range("A10").Value = 100
$ws.Range("A10").HorizontalAlignment = -4131

# 3) Delete the now-redundant "5k" / R11 row (row 19) entirely.
$ws.Rows.Item(19).Delete()

# 4) Match the author's final selection (Excel moved the active cell up
#    by one row after the deletion).
$ws.Range("A20").Select()
